$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Before adding the new row, copy row 60's current ("last row") date style
# onto the new row 61 so it becomes the new "last row" with date-only format.
$ws.Range("A60").Copy()
$ws.Range("A61").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 60 is no longer the last row, so it must switch to the regular
# date-time style used by all other data rows. Copy that style from A59.
$ws.Range("A59").Copy()
$ws.Range("A60").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Append the new data row (61) values.
$ws.Cells.Item(61, 1).Value = 45801
$ws.Cells.Item(61, 2).Value = 254
$ws.Cells.Item(61, 3).Value = 265
$ws.Cells.Item(61, 4).Value = 259
